$d = $word.ActiveDocument

# --- Locate the "By the help of microservices..." key-benefits list item ---
# (it is the paragraph that currently carries the _GoBack bookmark / is last
# in the "Key Benefits" numbered list)
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs.Item($i).Range.Text
    if ($ptxt -like "By the help of microservices*") {
        $targetIndex = $i
    }
}
if ($targetIndex -eq 0) {
    throw "Could not locate the 'By the help of microservices...' paragraph"
}

# --- Insert the first new list item right after it ---
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphAfter()
$idx1 = $targetIndex + 1
$d.Paragraphs.Item($idx1).Range.Text = "Once the whole application application is developed its easy to deployed them each other. And easy to identify the hot services and deploy them independently."

# --- Insert the second new list item right after the first one ---
$d.Paragraphs.Item($idx1).Range.InsertParagraphAfter()
$idx2 = $idx1 + 1
$d.Paragraphs.Item($idx2).Range.Text = "In the case of error the whole application doesn’t crash as they are not tightly coupled and only that service affects and when the error fixed the service start functioning from where it lefts."

# --- Append the trailing " " onto the end of the second new paragraph ---
$p2 = $d.Paragraphs.Item($idx2)
$endPos = $p2.Range.End - 1
$endR = $d.Range($endPos, $endPos)
$endR.InsertAfter(" ")

# --- Move the _GoBack bookmark so it sits right before that trailing space
#     (i.e. at the end of the real text, between the text run and the
#     space run) instead of at the end of the old paragraph ---
$p2b = $d.Paragraphs.Item($idx2)
$bmPos = $p2b.Range.End - 2
$bmR = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmR)
